$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7228
$ws.Range("K3").Value = 7505
$ws.Range("K4").Value = 1563
$ws.Range("K6").Value = 8286
$ws.Range("K7").Value = 25111

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 454
$ws.Range("K3").Value = 492
$ws.Range("K6").Value = 549
$ws.Range("K7").Value = 1634

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K6").Value = 338
$ws.Range("K7").Value = 1064

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 140
$ws.Range("K4").Value = 20
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 416

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 253
$ws.Range("K7").Value = 841

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 193
$ws.Range("K3").Value = 142
$ws.Range("K6").Value = 221
$ws.Range("K7").Value = 595

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 760
$ws.Range("K8").Value = 1634
$ws.Range("K10").Value = 146
$ws.Range("K15").Value = 257
$ws.Range("K18").Value = 168
$ws.Range("K19").Value = 732
$ws.Range("K20").Value = 615
$ws.Range("K26").Value = 33
$ws.Range("K27").Value = 239
$ws.Range("K29").Value = 1374
$ws.Range("K31").Value = 289
$ws.Range("K32").Value = 27
$ws.Range("K33").Value = 1064
$ws.Range("K34").Value = 142
$ws.Range("K36").Value = 319
$ws.Range("K37").Value = 841
$ws.Range("K39").Value = 31
$ws.Range("K42").Value = 926
$ws.Range("K47").Value = 169
$ws.Range("K52").Value = 649
$ws.Range("K54").Value = 490
$ws.Range("K57").Value = 98
$ws.Range("K59").Value = 44
$ws.Range("K60").Value = 144
$ws.Range("K65").Value = 595
$ws.Range("K72").Value = 123
$ws.Range("K73").Value = 225
$ws.Range("K77").Value = 169
$ws.Range("K78").Value = 305
$ws.Range("K85").Value = 1153
$ws.Range("K86").Value = 154
$ws.Range("K88").Value = 272
$ws.Range("K89").Value = 377
$ws.Range("K94").Value = 331
$ws.Range("K95").Value = 416
$ws.Range("K98").Value = 129
$ws.Range("K100").Value = 45
$ws.Range("K101").Value = 25111

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 88
$ws.Range("K7").Value = 289

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 356
$ws.Range("K6").Value = 274

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 76
$ws.Range("K3").Value = 111
$ws.Range("K7").Value = 490

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K6").Value = 399
$ws.Range("K7").Value = 1374

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 241
$ws.Range("K7").Value = 732

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 72
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 275
$ws.Range("K7").Value = 926

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 89
$ws.Range("K7").Value = 305

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 56

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 212
$ws.Range("K7").Value = 615

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 100
$ws.Range("K7").Value = 319

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 250
$ws.Range("K3").Value = 243
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 208
$ws.Range("K7").Value = 760

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 169

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 98
$ws.Range("K3").Value = 64
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 20
$ws.Range("K6").Value = 31

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 71
$ws.Range("K7").Value = 272

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value = 113
$ws.Range("K7").Value = 377

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 64
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 401
$ws.Range("K4").Value = 59
$ws.Range("K6").Value = 286
$ws.Range("K7").Value = 1153

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 70
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 169

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 238
$ws.Range("K7").Value = 649
